# Update the "F" column (想去人数 / interested count) values for specific
# events that were refreshed in the upstream data source. The same events
# appear both on the "展览" sheet and on the aggregated "全部类型" sheet,
# so the same value updates must be applied in both places.

$wb = $excel.ActiveWorkbook

# Map: sheet name -> list of (cell address, new value)
$updates = @{
    "展览" = @(
        @{ Cell = "F6";  Value = 344 },
        @{ Cell = "F8";  Value = 131 },
        @{ Cell = "F11"; Value = 5557 },
        @{ Cell = "F12"; Value = 38 },
        @{ Cell = "F13"; Value = 31 },
        @{ Cell = "F24"; Value = 1004 },
        @{ Cell = "F26"; Value = 1736 },
        @{ Cell = "F28"; Value = 31 }
    )
    "全部类型" = @(
        @{ Cell = "F8";  Value = 344 },
        @{ Cell = "F10"; Value = 131 },
        @{ Cell = "F13"; Value = 5557 },
        @{ Cell = "F14"; Value = 38 },
        @{ Cell = "F15"; Value = 31 },
        @{ Cell = "F34"; Value = 1004 },
        @{ Cell = "F36"; Value = 1736 },
        @{ Cell = "F38"; Value = 31 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range($entry.Cell).Value = $entry.Value
    }
}
